$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.325.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.31%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.67%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.72"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5153"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3965"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07834"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.115"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.374"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9971"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.356"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.807.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.58%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001080"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06580"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.58%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.027"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.325.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.76%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.219"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.467"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.89%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.017.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.87%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1099"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.068"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.14%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.592"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.54%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.653"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.63%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07184"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.204"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02357"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2185"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "

# Row 39
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.99%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.050"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.60%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6209"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.163"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.68%  "

# Row 43
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9993"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6006"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.304"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.741"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.926"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06802"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "
